$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (row 5), continuing the daily coverage-history log.
$ws.Range("A5").Value = 44057
$ws.Range("A5").NumberFormat = $ws.Range("A4").NumberFormat

$ws.Range("C5").Value = 161
$ws.Range("D5").Value = 73
$ws.Range("E5").Value = 1445
$ws.Range("F5").Value = 427
$ws.Range("G5").Value = 34
$ws.Range("H5").Value = 20
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 209
$ws.Range("K5").Value = 83

$ws.Range("M5").Formula = "=100*D5/C5"
$ws.Range("N5").Formula = "=100*F5/E5"
$ws.Range("O5").Formula = "=100*G5/C5"
$ws.Range("P5").Formula = "=100*I5/H5"
$ws.Range("Q5").Formula = "=100*K5/J5"

$ws.Range("M5:Q5").NumberFormat = $ws.Range("M4:Q4").NumberFormat

$ws.Calculate()

# Update selection to mirror the newly-entered row, as left by the author.
$ws.Range("B5:K5").Select()
